# Updates cryptos list (Price / Volume(1h) columns) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => Price (D), Volume(1h) (E)
$updates = @{
    2  = @{ D = "26.434.44"; E = "  +1.15%  " }
    3  = @{ D = "1.674.70";  E = "  +1.17%  " }
    5  = @{ D = "221.31";    E = "  +1.51%  " }
    6  = @{ D = "0.5353";    E = "  +1.00%  " }
    7  = @{ E = "  +0.77%  " }
    8  = @{ D = "0.2676";    E = "  +2.36%  " }
    9  = @{ D = "0.06415";   E = "  +1.29%  " }
    10 = @{ D = "21.03";     E = "  +2.97%  " }
    11 = @{ D = "0.07849";   E = "  +0.50%  " }
    12 = @{ D = "4.550";     E = "  +0.74%  " }
    13 = @{ D = "1.664.04";  E = "  -0.27%  " }
    14 = @{ D = "1.904.49";  E = "  +1.15%  " }
    15 = @{ D = "0.5665";    E = "  +3.22%  " }
    16 = @{ D = "0.0₅8213";  E = "  +0.06%  " }
    17 = @{ E = "  +1.69%  " }
    18 = @{ D = "26.465.22"; E = "  +1.26%  " }
    19 = @{ E = "  +0.77%  " }
    20 = @{ D = "4.734";     E = "  +2.85%  " }
    21 = @{ D = "199.02";    E = "  +4.05%  " }
    22 = @{ D = "10.37";     E = "  +2.88%  " }
    23 = @{ D = "6.086";     E = "  +1.09%  " }
    24 = @{ E = "  +0.75%  " }
    25 = @{ D = "146.77";    E = "  +1.05%  " }
    26 = @{ D = "0.1235";    E = "  +0.34%  " }
    27 = @{ D = "7.268";     E = "  +0.70%  " }
    28 = @{ D = "16.24";     E = "  +1.60%  " }
    29 = @{ D = "1.503";     E = "  +2.76%  " }
    30 = @{ D = "0.05900";   E = "  +2.31%  " }
    31 = @{ E = "  +0.89%  " }
    32 = @{ D = "3.590";     E = "  +0.80%  " }
    33 = @{ D = "3.319";     E = "  +1.56%  " }
    34 = @{ D = "1.623";     E = "  +1.50%  " }
    35 = @{ D = "0.9713";    E = "  +2.16%  " }
    36 = @{ D = "2.854";     E = "  +1.80%  " }
    37 = @{ D = "2.441";     E = "  +1.08%  " }
    38 = @{ D = "0.5837";    E = "  +1.48%  " }
    39 = @{ E = "  +0.61%  " }
    40 = @{ D = "1.082.02";  E = "  +3.66%  " }
    41 = @{ D = "5.931";     E = "  +2.34%  " }
    42 = @{ D = "0.8667";    E = "  +1.47%  " }
    43 = @{ E = "  +0.75%  " }
    44 = @{ D = "104.41";    E = "  -0.30%  " }
    45 = @{ D = "1.813.97";  E = "  +1.00%  " }
    46 = @{ D = "58.59";     E = "  +3.03%  " }
    47 = @{ E = "  -3.84%  " }
    48 = @{ E = "  +1.05%  " }
    49 = @{ D = "0.4417";    E = "  +1.64%  " }
    50 = @{ D = "8.058";     E = "  +2.31%  " }
    51 = @{ E = "  +0.38%  " }
}

foreach ($row in $updates.Keys) {
    $cells = $updates[$row]
    if ($cells.ContainsKey("D")) {
        # Force text format so numeric-looking price strings (e.g. "221.31")
        # remain text values instead of being auto-converted to numbers,
        # then restore General formatting so no visible format change remains.
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cells.D
        $cell.NumberFormat = "General"
    }
    if ($cells.ContainsKey("E")) {
        $ws.Range("E$row").Value = $cells.E
    }
}
